$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header column (BO) continuing the series after "06-sep"
$ws.Range("BO1").Value = "07-sep"

# New daily values for the added date column
$ws.Range("BO2").Value = 0
$ws.Range("BO3").Value = 20.559786839630078
$ws.Range("BO4").Value = 16.319772626934139
$ws.Range("BO5").Value = 18.888370164345623
$ws.Range("BO6").Value = 0
$ws.Range("BO7").Value = 17.456931160933685
$ws.Range("BO8").Value = 13.157271079161424
$ws.Range("BO9").Value = 15.241169988155876
$ws.Range("BO10").Value = 12.246121070953558
$ws.Range("BO11").Value = 13.690113140510944
$ws.Range("BO12").Value = 0
$ws.Range("BO13").Value = 8.2562461796624529
$ws.Range("BO14").Value = 0
$ws.Range("BO15").Value = 0
$ws.Range("BO16").Value = 13.374837060066978
$ws.Range("BO17").Value = 0
$ws.Range("BO18").Value = 0

# Match the author's final active selection cell
$ws.Range("BQ6").Select()
